$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 previously held "osnashennost_laboratornim_oborudovaniem" (B16).
# It is being replaced with the same text already used in row 15
# ("nalichie_informacii_o_vuze_v_internete"), which drops the now-unused
# shared string from the workbook.
$ws.Range("B16").Value = "nalichie_informacii_o_vuze_v_internete"

# Highlight row 12 (A12:B12) and row 16 (A16:B16) with a solid red fill.
# 255 == R=255,G=0,B=0 in Excel's BGR-packed Color long.
$ws.Range("A12:B12").Interior.Color = 255
$ws.Range("A16:B16").Interior.Color = 255

# Move/collapse the selection to A16.
$ws.Range("A16").Select()
